# Updated cryptos list values (Price / Volume(1h) columns, plus a two-row swap
# of the Avalanche / WrappedliquidstakedEther2.0 entries at rows 15-16).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of cell address -> new text value.
$updates = [ordered]@{
    'D2' = '59.862.44'
    'E2' = '  -1.72%  '
    'D3' = '2.302.55'
    'E3' = '  -2.90%  '
    'D4' = '1.00'
    'E4' = '  +0.02%  '
    'D5' = '540.74'
    'E5' = '  -1.54%  '
    'D6' = '128.75'
    'E6' = '  -3.54%  '
    'E7' = '  +0.05%  '
    'D8' = '0.569'
    'E8' = '  -3.43%  '
    'D9' = '2.300.68'
    'E9' = '  -2.89%  '
    'E10' = '  -1.10%  '
    'E11' = '  -0.71%  '
    'E12' = '  -0.87%  '
    'D13' = '0.331'
    'E13' = '  -1.82%  '
    'D14' = '59.780.25'
    'E14' = '  -1.69%  '
    'B15' = 'WrappedliquidstakedEther2.0'
    'C15' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D15' = '2.713.95'
    'E15' = '  -2.86%  '
    'B16' = 'Avalanche'
    'C16' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'D16' = '23.11'
    'E16' = '  -4.27%  '
    'E17' = '  -1.65%  '
    'D18' = '2.303.70'
    'E18' = '  -3.06%  '
    'D19' = '10.40'
    'E19' = '  -3.53%  '
    'D20' = '310.68'
    'E20' = '  -2.11%  '
    'E21' = '  -4.50%  '
    'D22' = '6.47'
    'E22' = '  -6.96%  '
    'E23' = '  -0.04%  '
    'D24' = '63.63'
    'E24' = '  +0.10%  '
    'D25' = '0.169'
    'E25' = '  -2.13%  '
    'E26' = '  -0.19%  '
    'D27' = '7.69'
    'E27' = '  -4.58%  '
    'E28' = '  -2.70%  '
    'D29' = '171.12'
    'E29' = '  -0.60%  '
    'D30' = '1.17'
    'E30' = '  +1.96%  '
    'E31' = '  -3.41%  '
    'D32' = '0.0₃0713'
    'E32' = '  -3.41%  '
    'D33' = '5.80'
    'E33' = '  -2.12%  '
    'E34' = '  -6.17%  '
    'D35' = '0.378'
    'E35' = '  -1.88%  '
    'E36' = '  -0.01%  '
    'D37' = '17.75'
    'E37' = '  -2.07%  '
    'D38' = '1.00'
    'E38' = '  +0.06%  '
    'D39' = '4.04'
    'E39' = '  -3.74%  '
    'D40' = '311.46'
    'E40' = '  -2.48%  '
    'D41' = '37.89'
    'E41' = '  -0.98%  '
    'D42' = '1.50'
    'E42' = '  -3.35%  '
    'D43' = '136.15'
    'E43' = '  -5.74%  '
    'E44' = '  -2.51%  '
    'D45' = '0.0935'
    'E45' = '  -2.42%  '
    'D46' = '0.565'
    'E46' = '  -0.29%  '
    'D47' = '18.45'
    'E47' = '  -4.95%  '
    'D48' = '0.0487'
    'E48' = '  -3.00%  '
    'D49' = '0.0212'
    'E49' = '  -1.65%  '
    'E50' = '  +1.76%  '
    'E51' = '  -0.52%  '
}

# These cells hold values that look numeric to Excel (e.g. '1.00', '10.40');
# without forcing a text format first, assigning them would silently coerce
# to a Double and drop the significant trailing zero (matches source: t="inlineStr").
$forceTextRefs = @('D4', 'D5', 'D6', 'D8', 'D13', 'D16', 'D19', 'D20', 'D22', 'D24', 'D25', 'D27', 'D29', 'D30', 'D33', 'D35', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D45', 'D46', 'D47', 'D48', 'D49')

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    if ($forceTextRefs -contains $ref) {
        $cell.NumberFormat = '@'
        $cell.Value = $updates[$ref]
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $updates[$ref]
    }
}
